$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brian")

# Row 3 held placeholder "Test" values in A3/B3. Finish adding this queue
# item the same way row 2 (Tag1 / foleyb25@gmail.com / Sender / Yes) was
# filled in: a tag name, the target address, the rule type, and whether it
# applies to attachments.
$ws.Range("A3").Value = "Tag2"
$ws.Range("B3").Value = "no-reply@revature.net"
$ws.Range("C3").Value = "Sender"
$ws.Range("D3").Value = "Yes"

# Mirror the mailto hyperlink that row 2's address cell (B2) already has.
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:no-reply@revature.net") | Out-Null
$ws.Range("B3").Style = "Hyperlink"

# Move the selection down to D4, past the row that was just completed.
$ws.Range("D4").Select() | Out-Null
